# Fix issue parsing openjml
# Updates the "all_tools" (aggregate) sheet and the "openjml" sheet with
# corrected num_snippets_warnings / num_warnings counts and the correlation
# statistics that were recomputed from them for rows 17-21 (the
# perc_correct_output / brain_deact_31 / brain_deact_32 / complexity_level /
# time_to_understand metrics). Also restores the "openjml" sheet's column I
# width that shifted as part of the same commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all_tools"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("all_tools")

$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 39
$ws.Range("I17").Value = -0.1647705109143269
$ws.Range("J17").Value = 0.4027546538976249
$ws.Range("K17").Value = -0.2341913484699036
$ws.Range("L17").Value = 0.3826574736700414

$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 39
$ws.Range("I18").Value = -0.1626978433639921
$ws.Range("J18").Value = 0.4043745903773081
$ws.Range("K18").Value = -0.2500726110965064
$ws.Range("L18").Value = 0.3502469232913418

$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 39
$ws.Range("K19").Value = -0.2741759952985793
$ws.Range("L19").Value = 0.3041353554760347

$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 39
$ws.Range("I20").Value = 0.3721433734379864
$ws.Range("J20").Value = 0.05729901570416417
$ws.Range("K20").Value = 0.443979491086295
$ws.Range("L20").Value = 0.08493122100226924

$ws.Range("F21").Value = 15
$ws.Range("G21").Value = 39
$ws.Range("I21").Value = -0.1265427670608828
$ws.Range("J21").Value = 0.5166373798159882
$ws.Range("K21").Value = -0.159684920338733
$ws.Range("L21").Value = 0.554695992273345

# ---------------------------------------------------------------------
# Sheet "openjml"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("openjml")

$ws2.Columns.Item(9).ColumnWidth = 19.8

$ws2.Range("F17").Value = 14
$ws2.Range("G17").Value = 30
$ws2.Range("I17").Value = 0.03774256780481986
$ws2.Range("J17").Value = 0.8508111952177798
$ws2.Range("K17").Value = 0.06715718224198458
$ws2.Range("L17").Value = 0.8048182163992386

$ws2.Range("F18").Value = 14
$ws2.Range("G18").Value = 30
$ws2.Range("I18").Value = -0.2608745973749755
$ws2.Range("J18").Value = 0.189459515366179
$ws2.Range("K18").Value = -0.3941471994340144
$ws2.Range("L18").Value = 0.1308813812945024

$ws2.Range("F19").Value = 14
$ws2.Range("G19").Value = 30
$ws2.Range("I19").Value = -0.3540440964374667
$ws2.Range("J19").Value = 0.07494263177455186
$ws2.Range("K19").Value = -0.4534975499279393
$ws2.Range("L19").Value = 0.07768807638224377

$ws2.Range("F20").Value = 14
$ws2.Range("G20").Value = 30
$ws2.Range("I20").Value = 0.14034022285596
$ws2.Range("J20").Value = 0.481699512487204
$ws2.Range("K20").Value = 0.1568611789958072
$ws2.Range("L20").Value = 0.5618076700772363

$ws2.Range("F21").Value = 14
$ws2.Range("G21").Value = 30
$ws2.Range("I21").Value = 0.07453559924999299
$ws2.Range("J21").Value = 0.7077285315990198
$ws2.Range("K21").Value = 0.07913380065856659
$ws2.Range("L21").Value = 0.770810224273788
